# Alterações Disp x Agravo (agora com Suspeita em Aberto e com Atendimento)
#
# - Adds two new header columns: "Situação" (G) and "Foi atendido" (H)
# - Backfills "Situação" / "Foi atendido" for the existing rows (2-15)
# - Reorders ids 1897696 / 1897693 (now row 8 / row 9 respectively)
# - Appends 5 new dispensation rows (16-20) for the new patients

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as plain literal text, bypassing
# Excel's "smart" date auto-detection (which would otherwise turn an
# unquoted ambiguous dd/mm/yyyy string like "11/05/2021" or "06/07/2021"
# into a real date serial). We do this by entering it as a formula that
# evaluates to the literal string, then collapsing the formula down to
# its cached value with a values-only paste - exactly like a user typing
# ="11/05/2021" and then Paste Special > Values over it.
function Set-LiteralText($cell, [string]$text) {
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# --- New header cells, copying the existing header formatting (bold,
#     border, centered) from the neighbouring F1 header cell ---
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)

$ws.Cells.Item(1, 7).Value = "Situação"
$ws.Cells.Item(1, 8).Value = "Foi atendido"

# --- Full data block (rows 2-20), columns A-H ---
# (B = data column whose string may be ambiguous enough to be mistaken
#  for a date; that's flagged so it gets the literal-text treatment.)
$rows = @(
    @(1872177, "19/04/2021", $false, 2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 82440,  "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1884539, "29/04/2021", $false, 2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 4143,   "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1897653, "11/05/2021", $true,  2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 55440,  "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1897686, "11/05/2021", $true,  2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 332561, "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1897689, "11/05/2021", $true,  2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 230722, "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1897691, "11/05/2021", $true,  2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 46858,  "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1897696, "11/05/2021", $true,  2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 78306,  "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1897693, "11/05/2021", $true,  2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 95249,  "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1897703, "11/05/2021", $true,  2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 85991,  "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1917669, "26/05/2021", $false, 2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 88102,  "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1955001, "21/06/2021", $false, 2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 86267,  "PACIENTE SEM SUSPEITA",            "SIM"),
    @(1957247, "23/06/2021", $false, 2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 246175, "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1958038, "23/06/2021", $false, 2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 16325,  "PACIENTE SEM SUSPEITA",            "NÃO"),
    @(1958036, "23/06/2021", $false, 2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 260075, "PACIENTE SEM SUSPEITA",            "SIM"),
    @(1964304, "29/06/2021", $false, 2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 226861, "PACIENTE COM SUSPEITA EM ABERTO",  "SIM"),
    @(1973147, "06/07/2021", $true,  2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 44582,  "PACIENTE COM SUSPEITA EM ABERTO",  "SIM"),
    @(1972855, "06/07/2021", $true,  2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 10335,  "PACIENTE COM SUSPEITA EM ABERTO",  "SIM"),
    @(1972871, "06/07/2021", $true,  2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 77902,  "PACIENTE COM SUSPEITA EM ABERTO",  "SIM"),
    @(1973400, "06/07/2021", $true,  2064081, "UBS DR MILTON BARONI DE BARRETOS", 1, 13337,  "PACIENTE COM SUSPEITA EM ABERTO",  "SIM")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]

    $dateCell = $ws.Cells.Item($r, 2)
    if ($row[2]) {
        Set-LiteralText $dateCell $row[1]
    } else {
        $dateCell.Value = $row[1]
    }

    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $r = $r + 1
}

$excel.CutCopyMode = 0
